$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edited transportation section: row 4 (GB/United Kingdom) -> replaced with DE/Germany, co2 updated to 0.35
$ws.Range("A4").Value = "DE"
$ws.Range("B4").Value = "Germany"
$ws.Range("C4").Value = 0.35

# Added diet section: new rows 5-7 duplicating country names with alternate (3-letter) country codes
$ws.Range("A5").Value = "IND"
$ws.Range("B5").Value = "India"
$ws.Range("C5").Value = 0.82

$ws.Range("A6").Value = "USA"
$ws.Range("B6").Value = "United States"
$ws.Range("C6").Value = 0.45

$ws.Range("A7").Value = "DEU"
$ws.Range("B7").Value = "Germany"
$ws.Range("C7").Value = 0.35

# Update selection to match the target view state
$ws.Range("F9").Select()
